$d = $word.ActiveDocument

function Add-JournalEntry($dateText, $bodyText) {
    # Append a new paragraph at the end of the document
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $endRange = $lastPara.Range
    $endRange.Collapse(0)
    $endRange.InsertParagraphAfter()

    # Locate the freshly created (empty) paragraph
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $insertPos = $newPara.Range.Start

    # Insert the bold "date" run
    $r1 = $d.Range($insertPos, $insertPos)
    $r1.InsertAfter($dateText)
    $boldRange = $d.Range($insertPos, $insertPos + $dateText.Length)
    $boldRange.Font.Bold = $true
    $boldRange.Font.Size = 12

    # Insert the regular "body" run right after it
    $bodyPos = $insertPos + $dateText.Length
    $r2 = $d.Range($bodyPos, $bodyPos)
    $r2.InsertAfter($bodyText)
    $normalRange = $d.Range($bodyPos, $bodyPos + $bodyText.Length)
    $normalRange.Font.Bold = $false
    $normalRange.Font.Size = 12
}

Add-JournalEntry "Mercredi 29 octobre 2014 : " "J’ai regardé le code effectuer par Xavier et je lui ai indiqué ce qui lui manquait pour que tout soit parfait et fonctionnel. "
Add-JournalEntry "Jeudi 30 octobre 2014 : " "J’ai vérifié le travail de tous mes coéquipiers et j’ai réglé les bugs présents dans le programme. "
